# Apply Chicago violent crime data update for 2025-12-20
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 6414
$ws.Range('L3').Value = 6908
$ws.Range('I4').Value = 1851
$ws.Range('L4').Value = 1715
$ws.Range('K5').Value = 593
$ws.Range('L6').Value = 5671
$ws.Range('I7').Value = 26321
$ws.Range('K7').Value = 27587
$ws.Range('L7').Value = 21114

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L2').Value = 71
$ws.Range('L7').Value = 238

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 429
$ws.Range('L3').Value = 493
$ws.Range('L4').Value = 95
$ws.Range('L6').Value = 338
$ws.Range('L7').Value = 1399

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L3').Value = 184
$ws.Range('L7').Value = 462

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L4').Value = 65
$ws.Range('L7').Value = 951

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L6').Value = 76
$ws.Range('L7').Value = 365

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L5').Value = 75
$ws.Range('L6').Value = 170
$ws.Range('L7').Value = 671
$ws.Range('L8').Value = 1399
$ws.Range('L10').Value = 140
$ws.Range('L11').Value = 350
$ws.Range('L16').Value = 48
$ws.Range('L19').Value = 576
$ws.Range('L20').Value = 533
$ws.Range('L23').Value = 221
$ws.Range('L29').Value = 1183
$ws.Range('L33').Value = 951
$ws.Range('L43').Value = 159
$ws.Range('L53').Value = 238
$ws.Range('L54').Value = 457
$ws.Range('L55').Value = 225
$ws.Range('L60').Value = 144
$ws.Range('I63').Value = 271
$ws.Range('K63').Value = 181
$ws.Range('L63').Value = 64
$ws.Range('L67').Value = 730
$ws.Range('L76').Value = 329
$ws.Range('L79').Value = 587
$ws.Range('L83').Value = 462
$ws.Range('L85').Value = 1047
$ws.Range('L90').Value = 224
$ws.Range('L96').Value = 234
$ws.Range('L99').Value = 365
$ws.Range('I101').Value = 26321
$ws.Range('K101').Value = 27587
$ws.Range('L101').Value = 21114

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L3').Value = 284
$ws.Range('L7').Value = 730

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L2').Value = 85
$ws.Range('L7').Value = 457

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 359
$ws.Range('L3').Value = 454
$ws.Range('L6').Value = 287
$ws.Range('L7').Value = 1183

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 208
$ws.Range('L3').Value = 176
$ws.Range('L7').Value = 576

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L4').Value = 43
$ws.Range('L7').Value = 329

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('L2').Value = 68
$ws.Range('L7').Value = 170

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('L4').Value = 5
$ws.Range('L7').Value = 140

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L3').Value = 74
$ws.Range('L7').Value = 225

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('L4').Value = 18
$ws.Range('L7').Value = 221

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L2').Value = 75
$ws.Range('L7').Value = 234

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L2').Value = 182
$ws.Range('L3').Value = 190
$ws.Range('L7').Value = 587

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 169
$ws.Range('L7').Value = 533

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 234
$ws.Range('L3').Value = 215
$ws.Range('L7').Value = 671

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L2').Value = 132
$ws.Range('L7').Value = 350

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('L3').Value = 21
$ws.Range('L7').Value = 75

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('L2').Value = 75
$ws.Range('L7').Value = 224

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('L4').Value = 9
$ws.Range('L7').Value = 144

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('L6').Value = 49
$ws.Range('L7').Value = 159

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L3').Value = 436
$ws.Range('L7').Value = 1047

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('L6').Value = 30
$ws.Range('L7').Value = 48
